# Updated cryptos list on Mon May 13 07:49:09 UTC 2024 with GitHub Actions
# Refreshes the Price (D) / Volume(1h) (E) columns for the crypto table,
# and restores the correct rank order for WrappedEther (row 18) and
# Polkadot (row 19), which had swapped positions.
#
# D-column values that are plain decimals (e.g. "595.34") are written with
# a leading apostrophe so Excel stores them as text (matching the sheet's
# existing text-typed Price column) instead of auto-converting them to
# numbers. Values already containing two dots (e.g. "61.695.39") or other
# non-numeric characters are left as-is since Excel can't parse them as
# numbers anyway.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.695.39"
$ws.Range("E2").Value = "  +1.52%  "

$ws.Range("D3").Value = "2.927.76"
$ws.Range("E3").Value = "  +0.65%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'595.34"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("D6").Value = "'141.87"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "2.927.77"
$ws.Range("E8").Value = "  +0.72%  "

$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("D10").Value = "'7.20"
$ws.Range("E10").Value = "  +4.50%  "

$ws.Range("E11").Value = "  +0.43%  "

$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("E13").Value = "  -0.81%  "

$ws.Range("D14").Value = "'32.99"
$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").Value = "3.413.16"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "61.523.46"
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "2.923.32"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'6.64"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").Value = "'433.42"
$ws.Range("E20").Value = "  +0.59%  "

$ws.Range("D21").Value = "'13.40"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").Value = "'0.668"
$ws.Range("E22").Value = "  -1.08%  "

$ws.Range("D23").Value = "'7.02"
$ws.Range("E23").Value = "  -0.69%  "

$ws.Range("D24").Value = "'81.08"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "'10.68"
$ws.Range("E25").Value = "  -1.53%  "

$ws.Range("D26").Value = "'2.12"
$ws.Range("E26").Value = "  -2.92%  "

$ws.Range("D27").Value = "'11.76"
$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").Value = "'2.10"
$ws.Range("E30").Value = "  -7.42%  "

$ws.Range("E31").Value = "  -1.78%  "

$ws.Range("D32").Value = "'26.19"

$ws.Range("E33").Value = "  +0.00%  "

$ws.Range("E34").Value = "  -3.03%  "

$ws.Range("D35").Value = "0.0₃0858"
$ws.Range("E35").Value = "  +0.27%  "

$ws.Range("E36").Value = "  -2.20%  "

$ws.Range("E37").Value = "  -1.00%  "

$ws.Range("D38").Value = "'49.22"
$ws.Range("E38").Value = "  -0.71%  "

$ws.Range("D39").Value = "'1.95"
$ws.Range("E39").Value = "  -0.80%  "

$ws.Range("E40").Value = "  -4.68%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("E42").Value = "  -2.56%  "

$ws.Range("D43").Value = "'0.274"
$ws.Range("E43").Value = "  -2.52%  "

$ws.Range("D44").Value = "'38.44"
$ws.Range("E44").Value = "  -6.86%  "

$ws.Range("D45").Value = "2.681.87"
$ws.Range("E45").Value = "  -0.41%  "

$ws.Range("D46").Value = "'133.35"
$ws.Range("E46").Value = "  -0.35%  "

$ws.Range("D47").Value = "'0.0336"
$ws.Range("E47").Value = "  -2.06%  "

$ws.Range("D48").Value = "'357.44"
$ws.Range("E48").Value = "  -4.66%  "

$ws.Range("D50").Value = "'22.71"
$ws.Range("E50").Value = "  -4.08%  "

$ws.Range("E51").Value = "  -2.11%  "
